# Updates Butte trap metadata workbook:
#  - rewrites several enumerated "Levels" definitions in the attribute sheet
#    to reflect new trap/site history (lower-cased levels, new site/position
#    options, extra visit/processed/trap-functioning codes)
#  - extends the historical min/max ranges for visitTime, rpmRevolutionsAtStart
#    and rpmRevolutionsAtEnd to cover the newly added historical data
#  - clears the stray alignment-only cell style that used to sit on B8

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("attribute")

# --- updated attribute_definition text (column B) -------------------------

$ws.Range("B4").Value = 'Name of the sampling site. Levels = c("parrot-phelan", "coleman national fish hatchery", "adams dam")'

$ws.Range("B5").Value = 'Position of trap within the site. Levels = c("canal trap box", "pp rst", "n/a", NA, "adams dam", "pp rst 2")'

$ws.Range("B7").Value = 'Work that was done during visit to trap. Levels = c("continue trapping", "start trap & begin trapping", "end trapping", "unplanned restart", "not applicable (n/a)", "service/adjust/clean trap", NA, "start trapping", "not recorded")'

$ws.Range("B8").Value = 'Describing if fish processed during trap visit. Levels = c("processed fish", "n/a; not a sampling visit", "no fish were caught",  "no catch data; fish released", NA, "not recorded")'

$ws.Range("B9").Value = 'Description of how well trap is functioning when visit to trap began. Levels = c("trap functioning normally", "trap functioning, but not normally", "trap not in service", "trap stopped functioning", "not recorded", NA)'

# --- widened historical min/max ranges -------------------------------------

# visitTime: earlier minimum, later maximum (older + newer trap visits)
$ws.Range("L6").Value = 35034.375
$ws.Range("M6").Value = 44929.489618055559

# rpmRevolutionsAtStart / rpmRevolutionsAtEnd maximum now much larger
$ws.Range("M12").Value = 3624
$ws.Range("M13").Value = 3624

# --- tidy up B8's one-off cell style ---------------------------------------
# B8 used to carry its own "applyFont + applyAlignment (no-op)" xf; re-assert
# the plain font-1/no-wrap formatting so it collapses back onto the same
# style as the rest of column A/B (xf index 4) instead of its own unique xf.

$ws.Range("B8").Font.Name = "Calibri"
$ws.Range("B8").Font.Size = 10
